$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Pass"
$ws.Range("C3").Value = "Pass"
$ws.Range("C4").Value = "Pass"
